$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose updated values look numeric ("231.59", "0.629", ...) need
# to be forced to Text storage so they stay inline/shared strings, just
# like every other cell in the Price/Volume columns (never real numbers).
# Apply a temporary Text format before writing, then restore the default
# "Normal" style so no visible formatting change is left behind.
$textCells = @('D5', 'D6', 'D7', 'D10', 'D11', 'D12', 'D17', 'D20', 'D21', 'D23', 'D25', 'D27', 'D28', 'D29', 'D32', 'D36', 'D38', 'D40', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated crypto market data (prices + 1h volume deltas),
# including the Cronos / Celestia / TrustWalletToken row re-ranking.
$ws.Range('D2').Value = '43.676.71'
$ws.Range('E2').Value = '  +5.21%  '
$ws.Range('D3').Value = '2.271.68'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '231.59'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').Value = '0.629'
$ws.Range('E6').Value = '  +1.82%  '
$ws.Range('D7').Value = '63.55'
$ws.Range('E7').Value = '  +6.47%  '
$ws.Range('E9').Value = '  +7.80%  '
$ws.Range('D10').Value = '0.105'
$ws.Range('E10').Value = '  +17.51%  '
$ws.Range('D11').Value = '57.34'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').Value = '25.90'
$ws.Range('E12').Value = '  +16.27%  '
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('D14').Value = '2.609.35'
$ws.Range('E14').Value = '  +2.78%  '
$ws.Range('E15').Value = '  +1.94%  '
$ws.Range('E16').Value = '  +5.19%  '
$ws.Range('D17').Value = '0.827'
$ws.Range('E17').Value = '  +4.25%  '
$ws.Range('D18').Value = '2.268.29'
$ws.Range('E18').Value = '  +2.56%  '
$ws.Range('D19').Value = '43.588.92'
$ws.Range('E19').Value = '  +4.99%  '
$ws.Range('D20').Value = '0.0000100'
$ws.Range('E20').Value = '  +11.42%  '
$ws.Range('D21').Value = '73.64'
$ws.Range('E21').Value = '  +2.35%  '
$ws.Range('E22').Value = '  +1.22%  '
$ws.Range('D23').Value = '250.18'
$ws.Range('E23').Value = '  +3.31%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').Value = '2.51'
$ws.Range('E25').Value = '  +6.92%  '
$ws.Range('E26').Value = '  +1.18%  '
$ws.Range('D27').Value = '9.91'
$ws.Range('E27').Value = '  +2.46%  '
$ws.Range('D28').Value = '172.52'
$ws.Range('E28').Value = '  +2.27%  '
$ws.Range('D29').Value = '20.98'
$ws.Range('E29').Value = '  +6.27%  '
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('E31').Value = '  +2.35%  '
$ws.Range('D32').Value = '2.79'
$ws.Range('E32').Value = '  +10.08%  '
$ws.Range('E33').Value = '  +1.36%  '
$ws.Range('E34').Value = '  +6.19%  '
$ws.Range('E35').Value = '  +2.47%  '
$ws.Range('D36').Value = '4.74'
$ws.Range('E36').Value = '  +2.34%  '
$ws.Range('E37').Value = '  +5.47%  '
$ws.Range('D38').Value = '3.81'
$ws.Range('E38').Value = '  +6.58%  '
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('D40').Value = '0.0248'
$ws.Range('E40').Value = '  +5.33%  '
$ws.Range('E41').Value = '  +0.26%  '
$ws.Range('D42').Value = '8.41'
$ws.Range('E42').Value = '  -1.65%  '
$ws.Range('D43').Value = '17.33'
$ws.Range('E43').Value = '  +6.00%  '
$ws.Range('B44').Value = 'Celestia'
$ws.Range('C44').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D44').Value = '10.40'
$ws.Range('E44').Value = '  +21.03%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '1.21'
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '0.0960'
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('D47').Value = '97.69'
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('D48').Value = '4.42'
$ws.Range('E48').Value = '  +4.20%  '
$ws.Range('D49').Value = '1.477.62'
$ws.Range('E49').Value = '  +0.99%  '
$ws.Range('D50').Value = '2.34'
$ws.Range('E50').Value = '  +4.78%  '
$ws.Range('E51').Value = '  +0.93%  '

# Restore the default cell style on the coerced cells now that the text
# value has stuck, so no stray formatting delta remains.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
